# Add a new "Spain" worksheet to the workbook, based on a copy of the
# "Italy" worksheet, then update its content for the Spain / Zettler
# market test data (NGC-3145/T2045).

$wb = $excel.ActiveWorkbook

$italy = $wb.Worksheets.Item("Italy")

# Copy the Italy sheet to the very end of the workbook -> becomes "Spain".
$italy.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))

$spain = $wb.Worksheets.Item($wb.Worksheets.Count)
$spain.Name = "Spain"

# Update the model number for Spain / Zettler market.
$spain.Range("B4").Value = "NGC-3145/T2045"

# The Spain sheet's columns re-flowed (column B narrower, column D wider)
# and rows 2-4 grew to a two-line height once the data settled.
$spain.Columns.Item(2).ColumnWidth = 14.3
$spain.Columns.Item(4).ColumnWidth = 26
$spain.Rows.Item(2).RowHeight = 28.8
$spain.Rows.Item(3).RowHeight = 28.8
$spain.Rows.Item(4).RowHeight = 28.8
$spain.Rows.Item(5).RowHeight = 28.8

# Italy's own selection is reset to the full used range (no specific
# active cell), and it's no longer the "tabSelected" sheet.
$italy.Activate()
$italy.Range("A1:D21").Select()

# Make the new sheet the active / selected sheet, with H3 selected.
$spain.Activate()
$spain.Range("H3").Select()
